$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of data to append below the existing header row.
$rows = @(
    @{ A = 45611; B = 4042080; C = "Creación Anticipo 409112186/1";     D = 4824; E = "IMEDISA ARTES GRAFICAS, S.L.U."; F = 0;      G = 361.94; H = -2470.21; I = 0; J = "sandra";  K = "yes"; L = 438004824; M = "yes"; N = "N"; O = "C"; P = "409112186/1" },
    @{ A = 45616; B = 4042912; C = "Compensación Anticipo 409112186/1"; D = 4824; E = "IMEDISA ARTES GRAFICAS, S.L.U."; F = 105.9;  G = 0;      H = -721.59;  I = 0; J = "MIREA93"; K = "yes"; L = 438004824; M = "yes"; N = "N"; O = "C"; P = "409112186/1" },
    @{ A = 45621; B = 4043419; C = "Compensación Anticipo 409112186/1"; D = 4824; E = "IMEDISA ARTES GRAFICAS, S.L.U."; F = 256.04; G = 0;      H = -465.55;  I = 0; J = "M.Jose";  K = "yes"; L = 438004824; M = "yes"; N = "N"; O = "C"; P = "409112186/1" },
    @{ A = 45656; B = 4047697; C = "Creación Anticipo 409127725/1";     D = 4824; E = "IMEDISA ARTES GRAFICAS, S.L.U."; F = 0;      G = 738.84; H = -738.84;  I = 0; J = "sandra";  K = "yes"; L = 438004824; M = "yes"; N = "N"; O = "C"; P = "409127725/1" }
)

$r = 2
$isFirstDateCell = $true
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row.A
    if ($isFirstDateCell) {
        # Register the lowercase datetime format first (matches the format
        # used elsewhere for "F.Mov."), then switch to the uppercase form
        # that is actually applied to the new "F.Mov." column cells.
        $ws.Cells.Item($r, 1).NumberFormat = "yyyy-mm-dd h:mm:ss"
        $isFirstDateCell = $false
    }
    $ws.Cells.Item($r, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Cells.Item($r, 2).Value = $row.B
    $ws.Cells.Item($r, 3).Value = $row.C
    $ws.Cells.Item($r, 4).Value = $row.D
    $ws.Cells.Item($r, 5).Value = $row.E
    $ws.Cells.Item($r, 6).Value = $row.F
    $ws.Cells.Item($r, 7).Value = $row.G
    $ws.Cells.Item($r, 8).Value = $row.H
    $ws.Cells.Item($r, 9).Value = $row.I
    $ws.Cells.Item($r, 10).Value = $row.J
    $ws.Cells.Item($r, 11).Value = $row.K
    $ws.Cells.Item($r, 12).Value = $row.L
    $ws.Cells.Item($r, 13).Value = $row.M
    $ws.Cells.Item($r, 14).Value = $row.N
    $ws.Cells.Item($r, 15).Value = $row.O
    $ws.Cells.Item($r, 16).Value = $row.P
    $r = $r + 1
}
